# Updating filtered feeds from workflow
# Appends two new rows (34 and 35) to the "Filtered Feeds" sheet, mirroring
# the existing row 33 (same keyword/title pairing), but pointing at the
# genomeweb.com and 360dx.com "cancer" variants of the Boehringer Ingelheim
# article link, each with its own hyperlink.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$keywords = "CDx"
$title = "Thermo Fisher Nabs FDA Approval for NGS-Based CDx for Boehringer Ingelheim Lung Cancer Drug"

$newLinks = @(
    "https://www.genomeweb.com/cancer/thermo-fisher-nabs-fda-approval-ngs-based-cdx-boehringer-ingelheim-lung-cancer-drug",
    "https://www.360dx.com/cancer/thermo-fisher-nabs-fda-approval-ngs-based-cdx-boehringer-ingelheim-lung-cancer-drug"
)

$xlUp = -4162
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End($xlUp).Row
$startRow = $lastRow + 1

for ($i = 0; $i -lt $newLinks.Count; $i++) {
    $row = $startRow + $i
    $link = $newLinks[$i]

    $cellA = $ws.Cells.Item($row, 1)
    $cellB = $ws.Cells.Item($row, 2)
    $cellC = $ws.Cells.Item($row, 3)

    $cellB.Value2 = $keywords
    $cellC.Value2 = $title

    $ws.Hyperlinks.Add($cellA, $link) | Out-Null
    $cellA.Style = "Hyperlink"
}
